$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new headers: I0 (col I) and IF (col J), matching the style of
# the other header cells (e.g. H1) by copying H1's formatting over.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Find the last used data row (column H currently holds data through row 22).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 8).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
